$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A45").Value = "Leo Parisi "
$ws.Range("B45").Value = "Daniele Dalbosco | IMONTAGNA"
$ws.Range("C45").Value = "Leonardo Viola | SHARK ATTACK"
$ws.Range("D45").Value = "Geremia  Carollo | FC SAVIGNANO"
$ws.Range("E45").Value = "Luca Frasca | Clitoriders"
$ws.Range("F45").Value = "Moris Benedetti | Gli Introvabili"
